$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 78, shifting existing rows 78:98 down to 79:99
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly data point
$ws.Cells.Item(78, 1).Value2  = 9
$ws.Cells.Item(78, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(78, 3).Value2  = "Metropolitana"
$ws.Cells.Item(78, 4).Value2  = 44711
$ws.Cells.Item(78, 5).Value2  = 13
$ws.Cells.Item(78, 6).Value2  = 100114007
$ws.Cells.Item(78, 7).Value2  = "Jengibre"
$ws.Cells.Item(78, 8).Value2  = "Sin especificar"
$ws.Cells.Item(78, 9).Value2  = "Primera"
$ws.Cells.Item(78, 10).Value2 = 610
$ws.Cells.Item(78, 11).Value2 = 17000
$ws.Cells.Item(78, 12).Value2 = 18000
$ws.Cells.Item(78, 13).Value2 = 17500
$ws.Cells.Item(78, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(78, 15).Value2 = "Perú"
$ws.Cells.Item(78, 16).Value2 = 1346
$ws.Cells.Item(78, 17).Value2 = 13
$ws.Cells.Item(78, 18).Value2 = "Hortaliza"
